$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "impuestos" column header
$ws.Range("D1").Value = "impuestos"

# Update existing values in row 2
$ws.Range("A2").Value = 10.43
$ws.Range("B2").Value = 10

# Add new tax value in row 2
$ws.Range("D2").Value = 0
